$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.937.23'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '2.458.77'
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '579.86'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.01%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '166.28'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.18%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -2.47%  '
$ws.Range('D9').Value = '2.458.52'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('E10').Value = '  -4.13%  '
$ws.Range('E11').Value = '  -1.02%  '
$ws.Range('E12').Value = '  -3.72%  '
$ws.Range('E13').Value = '  -2.90%  '
$ws.Range('D14').Value = '2.902.81'
$ws.Range('E14').Value = '  -2.47%  '
$ws.Range('E15').Value = '  -3.53%  '
$ws.Range('D16').Value = '66.809.05'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('E17').Value = '  -4.82%  '
$ws.Range('D18').Value = '2.461.24'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.64'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.04%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '353.69'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.02'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '69.17'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.86%  '
$ws.Range('E25').Value = '  -7.61%  '
$ws.Range('E26').Value = '  -7.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.90'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -9.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.989'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.07%  '
$ws.Range('D29').Value = '2.580.64'
$ws.Range('E29').Value = '  -1.71%  '
$ws.Range('D30').Value = '0.0₃0895'
$ws.Range('E30').Value = '  -7.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '507.28'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.48%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.76'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -5.93%  '
$ws.Range('E33').Value = '  -5.27%  '
$ws.Range('E34').Value = '  -6.27%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '157.92'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('E37').Value = '  -8.59%  '
$ws.Range('E38').Value = '  -0.46%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.42'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.93%  '
$ws.Range('E40').Value = '  -6.34%  '
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('E42').Value = '  -6.62%  '
$ws.Range('E44').Value = '  -6.66%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '38.70'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('E46').Value = '  -7.16%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '141.11'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('E48').Value = '  -6.01%  '
$ws.Range('E49').Value = '  -6.27%  '
$ws.Range('D50').Value = '0.0₆0253'
$ws.Range('E50').Value = '  -6.94%  '
$ws.Range('E51').Value = '  -6.77%  '
